$d = $word.ActiveDocument

# "Admin Server Login" -> "Admin Authorize"
# Replace "Server Login" with "Authorize" (the two runs "Server " and
# "Login" collapse into a single run whose text becomes "Authorize").
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Server Login", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Authorize", 2)
